$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update April 2021 value (H5)
$ws.Range("H5").Value = 100

# Update the August label to reflect the new "through" date
$ws.Range("A9").Value = "August (through 08-23)"

# Update August row (row 9) figures for 2016-2021 (B9/2015 unchanged)
$ws.Range("C9").Value = 52
$ws.Range("D9").Value = 62
$ws.Range("E9").Value = 41
$ws.Range("F9").Value = 34
$ws.Range("G9").Value = 135
$ws.Range("H9").Value = 118

# Update Total row (row 10) figures for 2016-2021 (B10/2015 unchanged)
$ws.Range("C10").Value = 354
$ws.Range("D10").Value = 527
$ws.Range("E10").Value = 466
$ws.Range("F10").Value = 338
$ws.Range("G10").Value = 756
$ws.Range("H10").Value = 1033
